# Adds the new article "Surah Al Imran, 186 - 200" as row 41 (article 40)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = @'
Surah Al Imran, 186 - 200
'@

$tags = @'
Connecting to our faith, Meaning of other books of religions, The right God, Definition of a believer
'@

$body = @'
h1: What is it to believe?
p: A human mind needs connections, justifications and logical evaluation of happenings to conclude. All conclusions drawn inside brain relate to the information it has. At times the information becomes dull and boring which makes the mind go asleep. Other times it is bright and easy to grasp, our mind wakes up. The information, its chaining and evidences are important factors to reach to a solid conclusion. 
p: We fight over God. Whose God is it that is right. Which religion is it that really matters. Labelling Islam as the only perfect religion without feeding the right information, connections and meanings to our brains do not make us believers. To become a believer is not an easier task. It takes a toll on human body. Huge amount of information needs to be evaluated, sorted, filtered and connected to finally reach the point where one can claim himself a ‘believer’.
h3: Paradise is for believers
p: It is true. Only believers will enter the paradise. Only people who are closer to the right God and His prophets are worthy to taste the flowing streams, green gardens, blossoming flowers and eternal life. The people who finally managed to find the right God will be rewarded. 
h3: Which God is the right God
p: All religions are fighting over the God. We all have taken the responsibility to defend our forefathers’ religions. We believe it is really easier to ‘Believe’. We also find it comforting to not think about believing because it takes huge effort. 
p: Labelling myself a believer only because my community forced it on me, doesn’t really make me a believer. I am just another pawn in the game of chess where other pawns are attacking me. I surrendered to the label of society. I am not still not a believer who can convince myself; <b>‘Why am I a believer?’.</b>
p: If we do not really believe, how can we convince the world to believe in our God. Lack of understandings, confusions and doubts do not go along with believing. <b>We can not defend our God until we ourselves believe in our God. </b> Without believing we can not conclude <b>‘Our God is the right God’.</b>
h3: Feeling the belief
p: The belief in the existence of God is a complicated step. It makes you atheist first, then it makes you a muslim. The label given by your family binds you to one religion. Later it makes you question your own sanity that leads you to stop believing. You are bound to a label not the belief. The moment you challenge this label your life can go back to ruins. It has become hard to become a logical believer.
p: Feeling the belief is the supreme task assigned to us. It is about feeling the Creator’s presence around us. We might call ourselves Christian, Jew or a Muslim but the fact remains there, the belief of having God around us surpasses all material logics. It is about believing in his all 4 x books. Bible if acted upon renders its reader eligible for heaven. 
p: I can not decide about other’s dealings with the creator. Torah if read properly, acted upon will lead the Jew to heavens. I can not really decide how the Creator interacts with people from other religions. 
h3: Conclusion
p: A Tsunami is approaching us. It is going to drown us away. On Judgement day, we will not be looking at a Jew standing beside us. <b>We will be looking at the Tsunami of injustices, mistakenly committed by us back on the planet earth. </b>
p: The God binds us to connect with Him only. Anybody who is connected with Him is His friend. He is the sole creator of heaven and earth who needs attention. He needs our time. 
p: Our affairs with other religions matter only because there is a God. If all other religions submit to God, we can not really blame them. We are bound to embrace them.
p: The real problem is not <b>“we are followers of Quran that we start feeling extra ordinarily right”.</b> The real problem is <b class=lavendar>we do not believe in Quran and its Writer. </span>
quote: O you who have believed, persevere and endure and remain stationed and fear Allah that you may be successful.<br>- Surah Al Imran verse 200
p.note: Guys, if you like this project. Please follow this project's page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.
'@

$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 43870
$ws.Range("B41").NumberFormat = "d-mmm-yy"
$ws.Range("C41").Value = $title
$ws.Range("F41").Value = $tags
$ws.Range("D41").Value = $body
$ws.Range("E41").Value = "Qasim Ali"

$ws.Rows.Item(41).RowHeight = 409.6

$ws.Range("D41").Select()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1

"done"
